$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value of 45181 (2023-09-12)
# for every data row (2 through 232). Update it to 45182 (2023-09-13).
$ws.Range("C2:C232").Value = 45182
